$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "anakum"
$ws.Range("C2").Value = "Ajit1234"
$ws.Range("E2").Value = "miam"
